# "average with safety stocks"
# Rescale InventoryCosts / BackorderCosts / LostSale on the Productdata sheet
# (divide by 2500), and zero out the safety-stock contribution for the last
# three buckets (rows 9-11, products 1-4) on the ForcastedStandardDeviation
# sheet.

$wb = $excel.ActiveWorkbook

# --- Productdata sheet: columns D (InventoryCosts), F (BackorderCosts), I (LostSale) ---
$pd = $wb.Worksheets.Item("Productdata")

$pd.Range("D2").Value = 0.0016
$pd.Range("F2").Value = 0.016
$pd.Range("I2").Value = 0.16

$pd.Range("D3").Value = 0.0028
$pd.Range("F3").Value = 0.028
$pd.Range("I3").Value = 0.28

$pd.Range("D4").Value = 0.0024
$pd.Range("F4").Value = 0.024
$pd.Range("I4").Value = 0.24

$pd.Range("D5").Value = 0.0012
$pd.Range("F5").Value = 0.012
$pd.Range("I5").Value = 0.12

$pd.Range("D6").Value = 0.0012
$pd.Range("F6").Value = 0.012
$pd.Range("I6").Value = 0.12

$pd.Range("D7").Value = 0.0012
$pd.Range("F7").Value = 0.012
$pd.Range("I7").Value = 0.12

$pd.Range("D8").Value = 0.0008
$pd.Range("F8").Value = 0.008
$pd.Range("I8").Value = 0.08

$pd.Range("D9").Value = 0.0004
$pd.Range("F9").Value = 0.004
$pd.Range("I9").Value = 0.04

$pd.Range("D10").Value = 0.0004
$pd.Range("F10").Value = 0.004
$pd.Range("I10").Value = 0.04

$pd.Range("D11").Value = 0.0004
$pd.Range("F11").Value = 0.004
$pd.Range("I11").Value = 0.04

# --- ForcastedStandardDeviation sheet: rows 9-11, columns B-E -> 0 ---
$sd = $wb.Worksheets.Item("ForcastedStandardDeviation")
$sd.Range("B9:E11").Value = 0
